# Updated cryptos list on Fri Jun 16 03:36:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the source's inline-string
# cell type) so numeric-looking strings like "1.000" or "25.524.04" are not
# auto-converted to numbers by Excel's smart cell entry.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") "25.524.04"
Set-TextValue $ws.Range("E2") "  +1.80%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.664.22"
Set-TextValue $ws.Range("E3") "  +0.80%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.9989"
Set-TextValue $ws.Range("E4") "  -0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "237.49"
Set-TextValue $ws.Range("E5") "  -0.28%  "

# Row 6
Set-TextValue $ws.Range("E6") "  +0.00%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.01%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.2628"

# Row 9
Set-TextValue $ws.Range("E9") "  +2.48%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07109"
Set-TextValue $ws.Range("E10") "  -1.16%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.659.98"
Set-TextValue $ws.Range("E11") "  +0.61%  "

# Row 12
Set-TextValue $ws.Range("D12") "14.78"
Set-TextValue $ws.Range("E12") "  -0.66%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.5882"
Set-TextValue $ws.Range("E13") "  -5.54%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.375"
Set-TextValue $ws.Range("E14") "  -4.83%  "

# Row 15
Set-TextValue $ws.Range("D15") "75.14"
Set-TextValue $ws.Range("E15") "  +2.49%  "

# Row 16
Set-TextValue $ws.Range("D16") "1.000"
Set-TextValue $ws.Range("E16") "  -0.02%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.000"
Set-TextValue $ws.Range("E17") "  +0.08%  "

# Row 18
Set-TextValue $ws.Range("D18") "25.511.87"
Set-TextValue $ws.Range("E18") "  +1.83%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.000006743"
Set-TextValue $ws.Range("E19") "  +1.80%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -0.18%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.875.33"
Set-TextValue $ws.Range("E21") "  +0.71%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.437"
Set-TextValue $ws.Range("E22") "  -2.20%  "

# Row 23
Set-TextValue $ws.Range("D23") "8.705"
Set-TextValue $ws.Range("E23") "  +0.89%  "

# Row 24
Set-TextValue $ws.Range("D24") "5.286"

# Row 25
Set-TextValue $ws.Range("D25") "134.81"
Set-TextValue $ws.Range("E25") "  +2.35%  "

# Row 26
Set-TextValue $ws.Range("D26") "15.04"
Set-TextValue $ws.Range("E26") "  +0.26%  "

# Row 27
Set-TextValue $ws.Range("E27") "  -1.10%  "

# Row 28
Set-TextValue $ws.Range("D28") "105.14"
Set-TextValue $ws.Range("E28") "  +1.63%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.708"
Set-TextValue $ws.Range("E29") "  +1.41%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +4.68%  "

# Row 31
Set-TextValue $ws.Range("D31") "3.658"
Set-TextValue $ws.Range("E31") "  +1.55%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.07709"
Set-TextValue $ws.Range("E32") "  -2.48%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.9992"
Set-TextValue $ws.Range("E33") "  -0.04%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.04249"
Set-TextValue $ws.Range("E34") "  -7.72%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.602"
Set-TextValue $ws.Range("E35") "  +0.31%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.6114"
Set-TextValue $ws.Range("E36") "  +6.21%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.9494"
Set-TextValue $ws.Range("E37") "  +0.51%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.602"
Set-TextValue $ws.Range("E38") "  -0.08%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.8642"
Set-TextValue $ws.Range("E39") "  +5.67%  "

# Row 40
Set-TextValue $ws.Range("E40") "  -0.04%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.853"
Set-TextValue $ws.Range("E41") "  +1.28%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.01466"
Set-TextValue $ws.Range("E42") "  -5.98%  "

# Row 43
Set-TextValue $ws.Range("D43") "96.81"
Set-TextValue $ws.Range("E43") "  -1.74%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.3765"
Set-TextValue $ws.Range("E44") "  +0.85%  "

# Row 45
Set-TextValue $ws.Range("D45") "4.831"
Set-TextValue $ws.Range("E45") "  +0.63%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.1122"
Set-TextValue $ws.Range("E46") "  -2.04%  "

# Row 47
Set-TextValue $ws.Range("D47") "6.203"
Set-TextValue $ws.Range("E47") "  +1.43%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +1.45%  "

# Row 49
Set-TextValue $ws.Range("D49") "29.74"
Set-TextValue $ws.Range("E49") "  -0.19%  "

# Row 50
Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.356"
Set-TextValue $ws.Range("E50") "  +2.16%  "

# Row 51
Set-TextValue $ws.Range("B51") "TrueUSD"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
Set-TextValue $ws.Range("D51") "1.001"
Set-TextValue $ws.Range("E51") "  -0.05%  "
